$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements (1-based row indices)
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "1488"
$t.Cell(6, 1).Range.Text  = "0.00127"
$t.Cell(7, 1).Range.Text  = "0.00019"
$t.Cell(8, 1).Range.Text  = "0.00006"
$t.Cell(9, 1).Range.Text  = "0.00026"
$t.Cell(10, 1).Range.Text = "0.00029"
$t.Cell(11, 1).Range.Text = "0.00040"
$t.Cell(12, 1).Range.Text = "0.27961"

# Rows that previously held 10 tab-separated values now hold a single value
$t.Cell(44, 1).Range.Text = "99.96"
$t.Cell(45, 1).Range.Text = "0.28"
$t.Cell(46, 1).Range.Text = "718"
